# @Xantam Update to Map032 scene update
# Add English translation column (C) alongside the existing Japanese
# dialogue column (A) for three newly-translated lines.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = "n\<n[3]>`"I'm a little hungry.
It's okay to drink your semen, if it's just a little.♥`""

$ws.Range("C12").Value = "n\<n[3]>Put it between my boobs
I'll make you all slimy ♥
Come on♥"

$ws.Range("C13").Value = "n\<n[3]>schlorp♥ schlorp♥ Come on, keep moving,
 I'm making your cock hard, right?
schlorp♥"
